# Auto-generated edit script applying the scheduled-runner market data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 13).Value = $null

$ws.Cells.Item(64, 8).Value = 4102.636
$ws.Cells.Item(64, 10).Value = 3848.4285
$ws.Cells.Item(64, 12).Value = 3848.4285
$ws.Cells.Item(64, 14).Value = -4344.4285

$ws.Cells.Item(67, 8).Value = 4102.636
$ws.Cells.Item(67, 10).Value = 3848.4285
$ws.Cells.Item(67, 12).Value = 3848.4285
$ws.Cells.Item(67, 14).Value = -5564.4285

$ws.Cells.Item(137, 8).Value = 1385.826
$ws.Cells.Item(137, 9).Value = 1291.5625
$ws.Cells.Item(137, 10).Value = 1601.2858
$ws.Cells.Item(137, 11).Value = 3874.6875
$ws.Cells.Item(137, 12).Value = 4803.857400000001
$ws.Cells.Item(137, 13).Value = -1324.6875
$ws.Cells.Item(137, 14).Value = -9903.857400000001

$ws.Cells.Item(138, 8).Value = 1466.763
$ws.Cells.Item(138, 9).Value = 898.9
$ws.Cells.Item(138, 10).Value = 1865.2632
$ws.Cells.Item(138, 11).Value = 2696.7
$ws.Cells.Item(138, 12).Value = 5595.7896
$ws.Cells.Item(138, 13).Value = 2443.3
$ws.Cells.Item(138, 14).Value = -15875.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 329.66666
$ws.Cells.Item(97, 9).Value = 329.66666
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 329.66666
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 166.33334
$ws.Cells.Item(97, 14).Value = $null

$ws.Cells.Item(102, 8).Value = 83333840
$ws.Cells.Item(102, 9).Value = 83333840
$ws.Cells.Item(102, 11).Value = 83333840
$ws.Cells.Item(102, 13).Value = -83332218

$ws.Cells.Item(132, 8).Value = 3112.9473
$ws.Cells.Item(132, 9).Value = 3379.375
$ws.Cells.Item(132, 10).Value = 2919.182
$ws.Cells.Item(132, 11).Value = 10138.125
$ws.Cells.Item(132, 12).Value = 8757.545999999998
$ws.Cells.Item(132, 13).Value = -7608.125
$ws.Cells.Item(132, 14).Value = -13817.546

$ws.Cells.Item(133, 8).Value = 27782.73
$ws.Cells.Item(133, 9).Value = 28000
$ws.Cells.Item(133, 10).Value = 27774.04
$ws.Cells.Item(133, 11).Value = 28000
$ws.Cells.Item(133, 12).Value = 27774.04
$ws.Cells.Item(133, 13).Value = -25470
$ws.Cells.Item(133, 14).Value = -32834.04

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3054.75
$ws.Cells.Item(86, 9).Value = 3140.1724
$ws.Cells.Item(86, 10).Value = 2700.8572
$ws.Cells.Item(86, 11).Value = 3140.1724
$ws.Cells.Item(86, 12).Value = 2700.8572
$ws.Cells.Item(86, 13).Value = -2017.1724
$ws.Cells.Item(86, 14).Value = -4946.8572

$ws.Cells.Item(89, 8).Value = 3054.75
$ws.Cells.Item(89, 9).Value = 3140.1724
$ws.Cells.Item(89, 10).Value = 2700.8572
$ws.Cells.Item(89, 11).Value = 15700.862
$ws.Cells.Item(89, 12).Value = 13504.286
$ws.Cells.Item(89, 13).Value = -10084.862
$ws.Cells.Item(89, 14).Value = -24736.286

$ws.Cells.Item(94, 8).Value = 35715704
$ws.Cells.Item(94, 9).Value = 41667990
$ws.Cells.Item(94, 10).Value = 1980
$ws.Cells.Item(94, 11).Value = 41667990
$ws.Cells.Item(94, 12).Value = 1980
$ws.Cells.Item(94, 13).Value = -41667539
$ws.Cells.Item(94, 14).Value = -2882

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1643
$ws.Cells.Item(31, 9).Value = 1464.5
$ws.Cells.Item(31, 11).Value = 1464.5
$ws.Cells.Item(31, 13).Value = -1169.5

$ws.Cells.Item(34, 8).Value = 1643
$ws.Cells.Item(34, 9).Value = 1464.5
$ws.Cells.Item(34, 11).Value = 1464.5
$ws.Cells.Item(34, 13).Value = -1262.5

$ws.Cells.Item(58, 8).Value = 928.6667
$ws.Cells.Item(58, 9).Value = 699.3929000000001
$ws.Cells.Item(58, 10).Value = 1731.125
$ws.Cells.Item(58, 11).Value = 699.3929000000001
$ws.Cells.Item(58, 12).Value = 1731.125
$ws.Cells.Item(58, 13).Value = -496.3929000000001
$ws.Cells.Item(58, 14).Value = -2137.125

$ws.Cells.Item(134, 8).Value = 25642970
$ws.Cells.Item(134, 9).Value = 30304964
$ws.Cells.Item(134, 11).Value = 90914892
$ws.Cells.Item(134, 13).Value = -90912357

$ws.Cells.Item(136, 8).Value = 928.6667
$ws.Cells.Item(136, 9).Value = 699.3929000000001
$ws.Cells.Item(136, 10).Value = 1731.125
$ws.Cells.Item(136, 11).Value = 2098.1787
$ws.Cells.Item(136, 12).Value = 5193.375
$ws.Cells.Item(136, 13).Value = 451.8212999999996
$ws.Cells.Item(136, 14).Value = -10293.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1464.0476
$ws.Cells.Item(5, 9).Value = 1464.0476
$ws.Cells.Item(5, 11).Value = 4392.142800000001
$ws.Cells.Item(5, 13).Value = -4280.142800000001

$ws.Cells.Item(40, 8).Value = 156.25
$ws.Cells.Item(40, 9).Value = 137.5
$ws.Cells.Item(40, 10).Value = 175
$ws.Cells.Item(40, 11).Value = 550
$ws.Cells.Item(40, 12).Value = 700
$ws.Cells.Item(40, 13).Value = -481
$ws.Cells.Item(40, 14).Value = -838

$ws.Cells.Item(44, 8).Value = 3499.5
$ws.Cells.Item(44, 9).Value = 3499
$ws.Cells.Item(44, 11).Value = 10497
$ws.Cells.Item(44, 13).Value = -10099

$ws.Cells.Item(64, 8).Value = 4329.231
$ws.Cells.Item(64, 9).Value = 2593.3333
$ws.Cells.Item(64, 10).Value = 4850
$ws.Cells.Item(64, 11).Value = 7779.999899999999
$ws.Cells.Item(64, 12).Value = 14550
$ws.Cells.Item(64, 13).Value = -7509.999899999999
$ws.Cells.Item(64, 14).Value = -15090

$ws.Cells.Item(67, 8).Value = 4329.231
$ws.Cells.Item(67, 9).Value = 2593.3333
$ws.Cells.Item(67, 10).Value = 4850
$ws.Cells.Item(67, 11).Value = 7779.999899999999
$ws.Cells.Item(67, 12).Value = 14550
$ws.Cells.Item(67, 13).Value = -6843.999899999999
$ws.Cells.Item(67, 14).Value = -16422

$ws.Cells.Item(69, 8).Value = 2024.4762
$ws.Cells.Item(69, 9).Value = 550
$ws.Cells.Item(69, 10).Value = 2179.6843
$ws.Cells.Item(69, 11).Value = 1650
$ws.Cells.Item(69, 12).Value = 6539.0529
$ws.Cells.Item(69, 13).Value = -839
$ws.Cells.Item(69, 14).Value = -8161.0529

$ws.Cells.Item(72, 8).Value = 2024.4762
$ws.Cells.Item(72, 9).Value = 550
$ws.Cells.Item(72, 10).Value = 2179.6843
$ws.Cells.Item(72, 11).Value = 4950
$ws.Cells.Item(72, 12).Value = 19617.1587
$ws.Cells.Item(72, 13).Value = -894
$ws.Cells.Item(72, 14).Value = -27729.1587

$ws.Cells.Item(103, 8).Value = 10519.272
$ws.Cells.Item(103, 9).Value = 331.25
$ws.Cells.Item(103, 10).Value = 16341
$ws.Cells.Item(103, 11).Value = 993.75
$ws.Cells.Item(103, 12).Value = 49023
$ws.Cells.Item(103, 13).Value = -114.75
$ws.Cells.Item(103, 14).Value = -50781

$ws.Cells.Item(135, 8).Value = 1464.0476
$ws.Cells.Item(135, 9).Value = 1464.0476
$ws.Cells.Item(135, 11).Value = 13176.4284
$ws.Cells.Item(135, 13).Value = -10641.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 13).Value = $null

$ws.Cells.Item(122, 8).Value = 377577
$ws.Cells.Item(122, 9).Value = 4500
$ws.Cells.Item(122, 10).Value = 750654
$ws.Cells.Item(122, 11).Value = 13500
$ws.Cells.Item(122, 12).Value = 2251962
$ws.Cells.Item(122, 13).Value = -11050
$ws.Cells.Item(122, 14).Value = -2256862

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1826.4546
$ws.Cells.Item(22, 9).Value = 1773.875
$ws.Cells.Item(22, 10).Value = 1966.6666
$ws.Cells.Item(22, 11).Value = 1773.875
$ws.Cells.Item(22, 12).Value = 1966.6666
$ws.Cells.Item(22, 13).Value = -1478.875
$ws.Cells.Item(22, 14).Value = -2556.6666

$ws.Cells.Item(27, 8).Value = 1826.4546
$ws.Cells.Item(27, 9).Value = 1773.875
$ws.Cells.Item(27, 10).Value = 1966.6666
$ws.Cells.Item(27, 11).Value = 1773.875
$ws.Cells.Item(27, 12).Value = 1966.6666
$ws.Cells.Item(27, 13).Value = -1666.875
$ws.Cells.Item(27, 14).Value = -2180.6666

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 13).Value = $null

$ws.Cells.Item(132, 8).Value = 20779.885
$ws.Cells.Item(132, 9).Value = 1311.7587
$ws.Cells.Item(132, 10).Value = 45326.652
$ws.Cells.Item(132, 11).Value = 3935.2761
$ws.Cells.Item(132, 12).Value = 135979.956
$ws.Cells.Item(132, 13).Value = -1405.2761
$ws.Cells.Item(132, 14).Value = -141039.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 20015
$ws.Cells.Item(43, 10).Value = 20015
$ws.Cells.Item(43, 12).Value = 20015
$ws.Cells.Item(43, 14).Value = -20313

$ws.Cells.Item(123, 8).Value = 53200
$ws.Cells.Item(123, 10).Value = 53200
$ws.Cells.Item(123, 12).Value = 53200
$ws.Cells.Item(123, 14).Value = -63000

$ws.Cells.Item(125, 8).Value = 65745
$ws.Cells.Item(125, 10).Value = 65745
$ws.Cells.Item(125, 12).Value = 65745
$ws.Cells.Item(125, 14).Value = -75585
